$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '29.688.27'
Set-TextValue 'E2' '  -1.47%  '

Set-TextValue 'D3' '2.097.30'
Set-TextValue 'E3' '  -0.48%  '

Set-TextValue 'D4' '1.010'
Set-TextValue 'E4' '  +0.32%  '

Set-TextValue 'D5' '343.38'
Set-TextValue 'E5' '  -1.75%  '

Set-TextValue 'E6' '  +0.19%  '

Set-TextValue 'D7' '0.5187'
Set-TextValue 'E7' '  +0.31%  '

Set-TextValue 'D8' '0.4376'
Set-TextValue 'E8' '  -1.82%  '

Set-TextValue 'D9' '53.54'
Set-TextValue 'E9' '  +1.69%  '

Set-TextValue 'D10' '0.09227'
Set-TextValue 'E10' '  +2.97%  '

Set-TextValue 'D11' '1.167'
Set-TextValue 'E11' '  -0.67%  '

Set-TextValue 'D12' '24.60'
Set-TextValue 'E12' '  -4.39%  '

Set-TextValue 'D13' '6.781'
Set-TextValue 'E13' '  +0.70%  '

Set-TextValue 'D14' '2.055.30'
Set-TextValue 'E14' '  -2.66%  '

Set-TextValue 'D15' '8.155'
Set-TextValue 'E15' '  -0.76%  '

Set-TextValue 'D16' '103.24'
Set-TextValue 'E16' '  +4.16%  '

Set-TextValue 'D17' '0.00001153'
Set-TextValue 'E17' '  +0.48%  '

Set-TextValue 'E18' '  +0.17%  '

Set-TextValue 'D19' '20.96'
Set-TextValue 'E19' '  +0.75%  '

Set-TextValue 'D20' '0.06671'
Set-TextValue 'E20' '  -0.42%  '

Set-TextValue 'D21' '1.008'
Set-TextValue 'E21' '  +0.30%  '

Set-TextValue 'D22' '6.211'
Set-TextValue 'E22' '  -0.54%  '

Set-TextValue 'D23' '29.723.26'
Set-TextValue 'E23' '  -1.72%  '

Set-TextValue 'D24' '12.55'
Set-TextValue 'E24' '  -2.44%  '

Set-TextValue 'E25' '  -1.67%  '

Set-TextValue 'D26' '2.323.83'
Set-TextValue 'E26' '  -1.44%  '

Set-TextValue 'D27' '21.90'
Set-TextValue 'E27' '  -0.28%  '

Set-TextValue 'D28' '161.92'
Set-TextValue 'E28' '  -0.38%  '

Set-TextValue 'D29' '2.492'
Set-TextValue 'E29' '  -1.84%  '

Set-TextValue 'D30' '133.59'
Set-TextValue 'E30' '  -0.04%  '

Set-TextValue 'D31' '1.128'
Set-TextValue 'E31' '  -4.01%  '

Set-TextValue 'D32' '1.687'
Set-TextValue 'E32' '  +3.50%  '

Set-TextValue 'D34' '6.197'
Set-TextValue 'E34' '  -0.98%  '

Set-TextValue 'D35' '3.954'

Set-TextValue 'D36' '6.355'
Set-TextValue 'E36' '  +7.53%  '

Set-TextValue 'D37' '10.42'
Set-TextValue 'E37' '  -0.12%  '

Set-TextValue 'D38' '0.02576'
Set-TextValue 'E38' '  -0.13%  '

Set-TextValue 'D39' '0.06717'
Set-TextValue 'E39' '  -1.76%  '

Set-TextValue 'D40' '0.6993'
Set-TextValue 'E40' '  +2.42%  '

Set-TextValue 'D41' '12.48'
Set-TextValue 'E41' '  -1.27%  '

Set-TextValue 'D42' '1.323'
Set-TextValue 'E42' '  +3.20%  '

Set-TextValue 'D43' '0.2217'
Set-TextValue 'E43' '  -4.00%  '

Set-TextValue 'D44' '0.6789'
Set-TextValue 'E44' '  +6.31%  '

Set-TextValue 'D45' '14.25'
Set-TextValue 'E45' '  -0.34%  '

Set-TextValue 'D46' '2.324'
Set-TextValue 'E46' '  +0.30%  '

Set-TextValue 'E47' '  -1.74%  '

Set-TextValue 'D48' '3.625'
Set-TextValue 'E48' '  -0.94%  '

Set-TextValue 'B49' 'EOS'
Set-TextValue 'C49' 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
Set-TextValue 'D49' '1.217'
Set-TextValue 'E49' '  -0.37%  '

Set-TextValue 'B50' 'WEMIXTOKEN'
Set-TextValue 'C50' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D50' '1.201'
Set-TextValue 'E50' '  +2.84%  '

Set-TextValue 'D51' '81.28'
Set-TextValue 'E51' '  -1.85%  '
